$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Row 2 - rent2
$ws.Range("A2").Value = "rent2"
$ws.Range("B2").Value = 900
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = Get-Date -Year 2025 -Month 7 -Day 19 -Hour 5 -Minute 30 -Second 10

# Row 3 - car
$ws.Range("A3").Value = "car"
$ws.Range("B3").Value = 45000
$ws.Range("C3").NumberFormat = "mm-dd-yy"
$ws.Range("C3").Value = Get-Date -Year 2025 -Month 7 -Day 19 -Hour 5 -Minute 30 -Second 10

# Row 4 - Pet food
$ws.Range("A4").Value = "Pet food"
$ws.Range("B4").Value = 123
$ws.Range("C4").NumberFormat = "mm-dd-yy"
$ws.Range("C4").Value = Get-Date -Year 2025 -Month 7 -Day 11 -Hour 5 -Minute 30 -Second 10

# Row 5 - Sister Wedding
$ws.Range("A5").Value = "Sister Wedding"
$ws.Range("B5").Value = 50000
$ws.Range("C5").NumberFormat = "mm-dd-yy"
$ws.Range("C5").Value = Get-Date -Year 2025 -Month 7 -Day 7 -Hour 5 -Minute 30 -Second 10
